# Apply the "feature count per condition" suffix to the Condicao column (C),
# refresh the AutoFilter over the full data range, update the hidden
# _FilterDatabase defined name to match, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each Condicao label to the number of features used for that condition.
$featureCounts = @{
    "Dataset Completo Balanceado - Nao Normalizado"     = 35
    "Dataset Completo Balanceado - Normalizado"         = 35
    "Dataset Completo Desbalanceado - Nao Normalizado"  = 35
    "Dataset Completo Desbalanceado - Normalizado"      = 35
    "PCA Balanceado"                                    = 12
    "PCA Desbalanceado"                                 = 12
    "ReliefF Balanceado - Nao Normalizado"              = 10
    "ReliefF Balanceado - Normalizado"                  = 10
    "ReliefF Desbalanceado - Nao Normalizado"           = 10
    "ReliefF Desbalanceado - Normalizado"               = 10
    "Branch and Bound Desbalanceado - Nao Normalizado"  = 6
    "Branch and Bound Balanceado - Nao Normalizado"     = 6
    "Branch and Bound Desbalanceado - Normalizado"      = 6
    "Branch and Bound Balanceado - Normalizado"         = 6
}

$lastRow = 169
$rng = $ws.Range("C2:C$lastRow")
$vals = $rng.Value2()

for ($i = 1; $i -le ($lastRow - 1); $i++) {
    $label = $vals[$i, 1]
    if ($featureCounts.ContainsKey($label)) {
        $count = $featureCounts[$label]
        $vals[$i, 1] = "$label ($count)"
    }
}

$rng.Value2 = $vals

# Re-apply the autofilter over the whole table (header + all data rows).
$fullRange = $ws.Range("A1:G$lastRow")
$fullRange.AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$$lastRow"
    }
}

# Move the active selection (no more scrolled-down topLeftCell override).
$ws.Range("C96").Select() | Out-Null
